# Update rubric scores on the "Basic Game rubric" sheet
# (powerup and jumpingEnemy functionality added -> higher scores)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basic Game rubric")

$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 2

# Move the active selection down to B4, matching the new state left by the author
$ws.Range("B4").Select()
